$wb = $excel.ActiveWorkbook

# --- Sheet "Means": update Total Cancer Risk (row 9) and Total Respiratory (row 10) ---
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 28
$wsMeans.Range("D9").Value = 17
$wsMeans.Range("E9").Value = 19
$wsMeans.Range("F9").Value = 19
$wsMeans.Range("G9").Value = 19

$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.3
$wsMeans.Range("D10").Value = 0.2
$wsMeans.Range("E10").Value = 0.2
$wsMeans.Range("F10").Value = 0.19
$wsMeans.Range("G10").Value = 0.18

# --- Sheet "Standard Deviations": update Total Cancer Risk (row 9) and Total Respiratory (row 10) ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 13
$wsSD.Range("D9").Value = 7.1
$wsSD.Range("E9").Value = 3
$wsSD.Range("F9").Value = 3
$wsSD.Range("G9").Value = 3.3

$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.11
$wsSD.Range("D10").Value = 0.000000000000000028
$wsSD.Range("E10").Value = 0.000000000000000012
$wsSD.Range("F10").Value = 0.03
$wsSD.Range("G10").Value = 0.039
